$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$values = @(
  @(1,1,"88-59="),
  @(1,2,"41-2="),
  @(1,3,"50-17="),
  @(1,4,"94-76="),
  @(1,5,"68+16="),
  @(2,1,"96-77="),
  @(2,2,"56-29="),
  @(2,3,"83-18="),
  @(2,4,"26+6="),
  @(2,5,"4+57="),
  @(3,1,"77-29="),
  @(3,2,"5+67="),
  @(3,3,"8+56="),
  @(3,4,"39+9="),
  @(3,5,"61-6="),
  @(4,1,"90-41="),
  @(4,2,"7+76="),
  @(4,3,"84-7="),
  @(4,4,"60-55="),
  @(4,5,"58+27="),
  @(5,1,"54+27="),
  @(5,2,"35+19="),
  @(5,3,"9+37="),
  @(5,4,"59+22="),
  @(5,5,"94-16="),
  @(6,1,"26+46="),
  @(6,2,"78+19="),
  @(6,3,"19+64="),
  @(6,4,"55+36="),
  @(6,5,"70-45="),
  @(7,1,"92-17="),
  @(7,2,"52+9="),
  @(7,3,"77+6="),
  @(7,4,"93-17="),
  @(7,5,"52-26="),
  @(8,1,"75+9="),
  @(8,2,"68+19="),
  @(8,3,"40-27="),
  @(8,4,"55-19="),
  @(8,5,"19+43="),
  @(9,1,"91-34="),
  @(9,2,"7+39="),
  @(9,3,"61-33="),
  @(9,4,"38+15="),
  @(9,5,"9+88="),
  @(10,1,"47+14="),
  @(10,2,"2+59="),
  @(10,3,"25+28="),
  @(10,4,"19+26="),
  @(10,5,"94-48="),
  @(11,1,"54-7="),
  @(11,2,"84-5="),
  @(11,3,"30-12="),
  @(11,4,"17+14="),
  @(11,5,"19+48="),
  @(12,1,"38+7="),
  @(12,2,"31-3="),
  @(12,3,"46-7="),
  @(12,4,"93-7="),
  @(12,5,"39+54="),
  @(13,1,"5+57="),
  @(13,2,"57+36="),
  @(13,3,"50-18="),
  @(13,4,"91-33="),
  @(13,5,"9+64="),
  @(14,1,"16+76="),
  @(14,2,"56+18="),
  @(14,3,"62-53="),
  @(14,4,"32-9="),
  @(14,5,"3+79="),
  @(15,1,"39+16="),
  @(15,2,"93-25="),
  @(15,3,"77+5="),
  @(15,4,"39+29="),
  @(15,5,"83-36="),
  @(16,1,"38+18="),
  @(16,2,"22+69="),
  @(16,3,"87+6="),
  @(16,4,"39+56="),
  @(16,5,"84-9="),
  @(17,1,"24+58="),
  @(17,2,"82-17="),
  @(17,3,"26+45="),
  @(17,4,"90-13="),
  @(17,5,"9+82="),
  @(18,1,"81-34="),
  @(18,2,"25+58="),
  @(18,3,"49+7="),
  @(18,4,"15+46="),
  @(18,5,"84-49="),
  @(19,1,"70-29="),
  @(19,2,"91-3="),
  @(19,3,"63-36="),
  @(19,4,"60-23="),
  @(19,5,"65-17="),
  @(20,1,"65+7="),
  @(20,2,"37+46="),
  @(20,3,"72-15="),
  @(20,4,"71-29="),
  @(20,5,"73-55="),
)

foreach ($item in $values) {
  $row = $item[0]
  $col = $item[1]
  $val = $item[2]
  $cell = $t.Cell($row, $col)
  $cell.Range.Text = $val
}

Write-Output "Done updating cells."